$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    # Force the cell to hold a literal text value even when the text looks
    # like a number (e.g. "-5.95"), mirroring how the source data is stored
    # as shared strings (t="s") rather than numeric cells in the workbook.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# NOTE: worksheet name lookup in this runtime is case-insensitive, and this
# workbook has sheets whose names differ only by case (e.g. "Vector_bf" vs
# "Vector_BF"). Use positional indices (matching xl/workbook.xml sheet
# order) instead of names to unambiguously address each sheet.
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---------------------------------------------------------------------
# Restricciones_del_follower: regenerated numeric example values
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "5.95 - y"
Set-TextValue $ws3.Range("B2") "-5.95"
Set-TextValue $ws3.Range("D2") "0.37"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "1.2"

$ws3.Range("A3").Value = "-0.6000000000000005 - x + y"
Set-TextValue $ws3.Range("B3") "-2.3999999999999995"
Set-TextValue $ws3.Range("D3") "0.44"
Set-TextValue $ws3.Range("E3") "6.2"
Set-TextValue $ws3.Range("F3") "7.0"

$ws3.Range("A4").Value = "-17.25 + x + 2y"
Set-TextValue $ws3.Range("B4") "5.25"
Set-TextValue $ws3.Range("D4") "0.0"
Set-TextValue $ws3.Range("F4") "4.4"

$ws3.Range("A5").Value = "-15.549999999999999 + 4x - y"
Set-TextValue $ws3.Range("B5") "3.4499999999999993"
Set-TextValue $ws3.Range("D5") "0.07"
Set-TextValue $ws3.Range("E5") "0"
Set-TextValue $ws3.Range("F5") "2.5"

# ---------------------------------------------------------------------
# Punto_modificado: regenerated (x, y) point
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "5.35"
Set-TextValue $ws4.Range("B2") "5.95"

# ---------------------------------------------------------------------
# Vector_bf: regenerated value
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-1.0"

# ---------------------------------------------------------------------
# Vector_BF: regenerated values
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-2.000000000000001"
Set-TextValue $ws6.Range("A3") "-21.6"
